$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.726.87'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '1.601.42'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.63'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.69'
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '1.826.06'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '1.605.44'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.523'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.11'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '0.0₃0738'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '210.48'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.17'
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.26'
$ws.Range("E22").Value = '  -2.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.99'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.67'
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.37'
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0510'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("E32").Value = '  +0.72%  '
$ws.Range("D33").Value = '1.290.38'
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.603'
$ws.Range("E36").Value = '  -2.97%  '
$ws.Range("E37").Value = '  +10.50%  '
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.76'
$ws.Range("E43").Value = '  -1.57%  '
$ws.Range("D44").Value = '1.738.08'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.56'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("E46").Value = '  -1.57%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0515'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.41'
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("E51").Value = '  +0.84%  '
